$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (column G) values, replacing the old Strike# values.
$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 0
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
